# Update manual status column (I) for the two rows that previously held
# the numeric value 8 (manualStatus) so that they instead hold the text
# label "[8]" - matching the new shared string entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 / Row 27: manualStatus column (I) becomes a text label "[8]"
$ws.Range("I26").Value = "[8]"
$ws.Range("I27").Value = "[8]"

# Widen the fastqFileName column (F) so the long file names are readable
$ws.Columns.Item(6).ColumnWidth = 67.1

# Tighten up row 27's height slightly
$ws.Rows.Item(27).RowHeight = 13.8

# Move the active selection to I27 to reflect where editing finished
$ws.Range("I27").Select()
